$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Patient data rows to append to Sheet1, matching the column headers already
# present in row 1 (id, bach_*, marcus_*, park_* scores).
$patientData = @(
    @(1, 50, 1, 30, 0, 0, 0, 0, 50, 1, 30, 1, 0, 0, 0, 50, 2, 2, 0, 0, 0),
    @(2, 65, 20, 50, 0, 0, 0, 0, 65, 1, 50, 20, 0, 0, 0, 65, 3, 3, 1, 1, 1)
)

for ($r = 0; $r -lt $patientData.Length; $r++) {
    $rowValues = $patientData[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $cell = $ws1.Cells.Item($r + 2, $c + 1)
        $cell.Value = $rowValues[$c]
        $cell.HorizontalAlignment = -4108
        $cell.Font.Color = 0
    }
}

# Select the newly written patient data range on Sheet1.
$ws1.Range("A2:U3").Select()

# Add a second worksheet (after Sheet1) that will hold the computed risk
# scores written out from the patient data above.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Patient Risk Scores"
